# Added late days to demo data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. "projects" sheet: add a new "C Project" column (E) mirroring
#    the m/p/x grade columns used elsewhere, with a few entries
#    carrying a "|<late days>" suffix.
# ---------------------------------------------------------------
$projects = $wb.Worksheets.Item("projects")

$projects.Range("E1").Value = "C Project"
$projects.Range("E2").Value = "c"
$projects.Range("E3").Value = "m"
$projects.Range("E4").Value = "m"
$projects.Range("E5").Value = "m"
$projects.Range("E6").Value = "m"
$projects.Range("E7").Value = "m|1"
$projects.Range("E8").Value = "m"
$projects.Range("E9").Value = "p"
$projects.Range("E10").Value = "p"
$projects.Range("E11").Value = "x"
$projects.Range("E12").Value = "."
$projects.Range("E13").Value = "."
$projects.Range("E14").Value = "x"
$projects.Range("E15").Value = "m"
$projects.Range("E16").Value = "p|5"
$projects.Range("E17").Value = "m"
$projects.Range("E18").Value = "m"
$projects.Range("E19").Value = "m|3"
$projects.Range("E20").Value = "m|4"
$projects.Range("E21").Value = "p"
$projects.Range("E22").Value = "x"
$projects.Range("E23").Value = "x"
$projects.Range("E24").Value = "x."
$projects.Range("E25").Value = "."
$projects.Range("E26").Value = "."
$projects.Range("E27").Value = "m"
$projects.Range("E28").Value = "p"

# ---------------------------------------------------------------
# 2. "learningObjectives" sheet: record late days on a handful of
#    previously-graded "Get/Post" (column F) entries.
# ---------------------------------------------------------------
$learningObjectives = $wb.Worksheets.Item("learningObjectives")

$learningObjectives.Range("F11").Value = "pm|3"
$learningObjectives.Range("F13").Value = "m|3"
$learningObjectives.Range("F16").Value = "m|2"
$learningObjectives.Range("F17").Value = "p|2"
$learningObjectives.Range("F19").Value = "pm|1"
$learningObjectives.Range("F20").Value = "m|2"

# ---------------------------------------------------------------
# 3. Cursor / active-tab bookkeeping to match the saved workbook
#    state: "learningObjectives" is left as the active sheet with
#    F21 selected, while "projects" is left with E29 selected.
# ---------------------------------------------------------------
$projects.Range("E29").Select()
$learningObjectives.Range("F21").Select()
$learningObjectives.Activate()
